$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated roster table (player name, position(s), team) matching the
# refreshed data from the source sheet. Row 19 is new; the rest are
# overwritten in place with the refreshed values.
$data = @(
    @("Fred VanVleet", "PG", "Houston Rockets"),
    @("Anfernee Simons", "PG,SG", "Portland Trail Blazers"),
    @("James Harden", "PG,SG", "LA Clippers"),
    @("Anthony Edwards", "SG,SF", "Minnesota Timberwolves"),
    @("Amen Thompson", "SG,SF", "Houston Rockets"),
    @("Jaren Jackson Jr.", "PF,C", "Memphis Grizzlies"),
    @("Bilal Coulibaly", "SG,SF", "Washington Wizards"),
    @("Giannis Antetokounmpo", "PF,C", "Milwaukee Bucks"),
    @("Zion Williamson", "PF,C", "New Orleans Pelicans"),
    @("Nicolas Claxton", "C", "Brooklyn Nets"),
    @("Jayson Tatum", "SF,PF", "Boston Celtics"),
    @("Keyonte George", "PG,SG", "Utah Jazz"),
    @("Ivica Zubac", "C", "LA Clippers"),
    @("Mason Plumlee", "C", "Phoenix Suns"),
    @("Paul George", "SG,SF,PF", "Philadelphia 76ers"),
    @("Jonathan Kuminga", "SF,PF", "Golden State Warriors"),
    @("Goga Bitadze", "C", "Orlando Magic"),
    @("Bradley Beal", "PG,SG,SF", "Phoenix Suns")
)

$row = 2
foreach ($record in $data) {
    $ws.Cells.Item($row, 1).Value = $record[0]
    $ws.Cells.Item($row, 2).Value = $record[1]
    $ws.Cells.Item($row, 3).Value = $record[2]
    $row++
}
